$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values are plain numeric-looking strings in the source data
# (inline strings, not real numbers). Force text format so COM does not silently
# coerce them to numbers, then restore the default "Normal" style so no stray
# cell-style index is introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "54.348.95"
$ws.Range("E2").Value = "  -2.69%  "

Set-TextValue $ws.Range("D3") "2.285.22"
$ws.Range("E3").Value = "  -3.09%  "

Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.20%  "

Set-TextValue $ws.Range("D5") "493.17"
$ws.Range("E5").Value = "  -2.23%  "

Set-TextValue $ws.Range("D6") "127.01"
$ws.Range("E6").Value = "  -2.49%  "

Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.92%  "

Set-TextValue $ws.Range("D9") "2.284.30"
$ws.Range("E9").Value = "  -3.64%  "

Set-TextValue $ws.Range("D10") "0.0944"
$ws.Range("E10").Value = "  -2.72%  "

$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("E13").Value = "  -3.66%  "

Set-TextValue $ws.Range("D14") "2.673.42"
$ws.Range("E14").Value = "  -3.74%  "

Set-TextValue $ws.Range("D15") "21.57"
$ws.Range("E15").Value = "  +0.27%  "

Set-TextValue $ws.Range("D16") "54.225.60"
$ws.Range("E16").Value = "  -2.83%  "

$ws.Range("E17").Value = "  -2.48%  "

Set-TextValue $ws.Range("D18") "2.272.16"
$ws.Range("E18").Value = "  -1.50%  "

Set-TextValue $ws.Range("D19") "9.98"
$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("E20").Value = "  +0.91%  "

Set-TextValue $ws.Range("D21") "303.43"
$ws.Range("E21").Value = "  -2.28%  "

$ws.Range("E22").Value = "  +3.88%  "

$ws.Range("E23").Value = "  +0.12%  "

Set-TextValue $ws.Range("D24") "5.34"
$ws.Range("E24").Value = "  -3.66%  "

$ws.Range("E25").Value = "  -2.83%  "

$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("E27").Value = "  +0.93%  "

Set-TextValue $ws.Range("D28") "2.395.11"

$ws.Range("E29").Value = "  +2.49%  "

Set-TextValue $ws.Range("D30") "7.07"
$ws.Range("E30").Value = "  -0.94%  "

Set-TextValue $ws.Range("D31") "169.01"
$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("E32").Value = "  -2.42%  "

$ws.Range("E33").Value = "  -3.50%  "

$ws.Range("E34").Value = "  +2.22%  "

Set-TextValue $ws.Range("D35") "0.998"
$ws.Range("E35").Value = "  -0.16%  "

Set-TextValue $ws.Range("D36") "0.999"
$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("E38").Value = "  -0.40%  "

$ws.Range("E39").Value = "  +2.32%  "

Set-TextValue $ws.Range("D40") "0.865"
$ws.Range("E40").Value = "  +3.08%  "

$ws.Range("E41").Value = "  -0.49%  "

Set-TextValue $ws.Range("D42") "35.47"
$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("E43").Value = "  +1.15%  "

$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "128.78"
$ws.Range("E45").Value = "  +2.25%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D46") "3.34"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("E47").Value = "  -1.57%  "

Set-TextValue $ws.Range("D48") "0.0892"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("E49").Value = "  -2.46%  "

Set-TextValue $ws.Range("D50") "239.10"
$ws.Range("E50").Value = "  -0.47%  "

Set-TextValue $ws.Range("D51") "0.0478"
$ws.Range("E51").Value = "  +0.02%  "
